$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 with new station_id based data
$ws.Range("A2").Value = "clbhmthk20000v4s4f4yu874a"
$ws.Range("B2").Value = "clbgnzizb0000v4ag550yepfe"
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = 10

# Update row 3 with new station_id based data
$ws.Range("A3").Value = "clbhmthk20000v4s4f4yu874a"
$ws.Range("B3").Value = "clbgnzizb0000v4ag550yepfe"
$ws.Range("C3").Value = 2
$ws.Range("D3").Value = 10

# Remove the now-obsolete rows 4 through 12
$ws.Range("A4:D12").ClearContents()

# Keep ignoredErrors definitions aligned with the shrunk data range
$ws.Range("A1:D3").Errors.Item(9).Ignore = $true
